# Update "想去人数" (people-interested count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - first worksheet
$wsExhibition = $wb.Worksheets.Item(1)
$wsExhibition.Range("F2").Value = 5345
$wsExhibition.Range("F6").Value = 808
$wsExhibition.Range("F7").Value = 318

# Sheet "演出" (Performance) - second worksheet
$wsPerformance = $wb.Worksheets.Item(2)
$wsPerformance.Range("F3").Value = 11

# Sheet "全部类型" (All types) - fourth worksheet
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 5345
$wsAll.Range("F6").Value = 808
$wsAll.Range("F8").Value = 318
$wsAll.Range("F10").Value = 11
